# The site was re-crawled; the product "Super Soft Premium 99% Water 4x 50ST"
# (row 4) disappeared from the listing between the two crawls, and every
# subsequent product shifted up by one row. All rows also get the new
# crawl timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for the product that dropped out of the listing; Excel
# shifts everything below it up by one (same as selecting row 4 and doing
# Delete > Entire Row).
$ws.Rows(4).Delete()

# Stamp every remaining data row (2 through the new last row, 38) with the
# updated crawl timestamp.
$ws.Range("O2:O38").Value = "2022-07-20 20:58:14"
